{"js": "// Replace each \"A\u00d7B=\" multiplication prompt with its updated numbers.\n// The mapping below preserves document order exactly as seen in the diff.\nconst replacements = [\n  [\"933\u00d74=\", \"425\u00d77=\"],\n  [\"738\u00d79=\", \"468\u00d79=\"],\n  [\"161\u00d74=\", \"766\u00d77=\"],\n  [\"518\u00d76=\", \"173\u00d78=\"],\n  [\"567\u00d79=\", \"101\u00d76=\"],\n  [\"841\u00d79=\", \"173\u00d74=\"],\n  [\"283\u00d77=\", \"562\u00d73=\"],\n  [\"359\u00d75=\", \"421\u00d75=\"],\n  [\"256\u00d74=\", \"605\u00d75=\"],\n  [\"171\u00d79=\", \"291\u00d75=\"],\n  [\"997\u00d75=\", \"252\u00d79=\"],\n  [\"600\u00d72=\", \"947\u00d76=\"],\n  [\"906\u00d78=\", \"897\u00d75=\"],\n  [\"281\u00d76=\", \"565\u00d79=\"],\n  [\"685\u00d79=\", \"694\u00d79=\"],\n  [\"378\u00d72=\", \"386\u00d76=\"],\n  [\"942\u00d74=\", \"655\u00d79=\"],\n  [\"955\u00d76=\", \"732\u00d79=\"],\n  [\"350\u00d79=\", \"655\u00d74=\"],\n  [\"913\u00d75=\", \"401\u00d73=\"],\n  [\"939\u00d74=\", \"709\u00d78=\"],\n  [\"224\u00d78=\", \"353\u00d78=\"],\n  [\"998\u00d75=\", \"638\u00d72=\"],\n  [\"366\u00d77=\", \"973\u00d79=\"],\n  [\"304\u00d76=\", \"632\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"A\u00d7B=\" multiplication prompt with its updated numbers.\n# The mapping below preserves document order exactly as seen in the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"933\u00d74=\", \"425\u00d77=\"),\n    @(\"738\u00d79=\", \"468\u00d79=\"),\n    @(\"161\u00d74=\", \"766\u00d77=\"),\n    @(\"518\u00d76=\", \"173\u00d78=\"),\n    @(\"567\u00d79=\", \"101\u00d76=\"),\n    @(\"841\u00d79=\", \"173\u00d74=\"),\n    @(\"283\u00d77=\", \"562\u00d73=\"),\n    @(\"359\u00d75=\", \"421\u00d75=\"),\n    @(\"256\u00d74=\", \"605\u00d75=\"),\n    @(\"171\u00d79=\", \"291\u00d75=\"),\n    @(\"997\u00d75=\", \"252\u00d79=\"),\n    @(\"600\u00d72=\", \"947\u00d76=\"),\n    @(\"906\u00d78=\", \"897\u00d75=\"),\n    @(\"281\u00d76=\", \"565\u00d79=\"),\n    @(\"685\u00d79=\", \"694\u00d79=\"),\n    @(\"378\u00d72=\", \"386\u00d76=\"),\n    @(\"942\u00d74=\", \"655\u00d79=\"),\n    @(\"955\u00d76=\", \"732\u00d79=\"),\n    @(\"350\u00d79=\", \"655\u00d74=\"),\n    @(\"913\u00d75=\", \"401\u00d73=\"),\n    @(\"939\u00d74=\", \"709\u00d78=\"),\n    @(\"224\u00d78=\", \"353\u00d78=\"),\n    @(\"998\u00d75=\", \"638\u00d72=\"),\n    @(\"366\u00d77=\", \"973\u00d79=\"),\n    @(\"304\u00d76=\", \"632\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
